$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-28 18:48:15"
$ws.Range("H2").Value = "'76%"
$ws.Range("H6").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("O2").Value = "3.0 °C"
$ws.Range("E3").Value = "2026-02-28 18:48:18"
$ws.Range("I3").Value = "0.2 mm"
$ws.Range("N3").Value = "-2.8 °C 18:29 TU"
$ws.Range("O3").Value = "-1.1 °C"
$ws.Range("E4").Value = "2026-02-28 18:48:20"
$ws.Range("K4").Value = "6.0 MJ/m2"
$ws.Range("E5").Value = "2026-02-28 18:48:23"
$ws.Range("N5").Value = "-2.6 °C 18:22 TU"
$ws.Range("O5").Value = "-1.2 °C"
$ws.Range("E6").Value = "2026-02-28 18:48:25"
$ws.Range("E7").Value = "2026-02-28 18:48:28"
$ws.Range("E8").Value = "2026-02-28 18:48:30"
$ws.Range("E9").Value = "2026-02-28 18:48:33"
$ws.Range("E10").Value = "2026-02-28 18:48:34"
$ws.Range("E11").Value = "2026-02-28 18:48:35"
$ws.Range("E12").Value = "2026-02-28 18:48:36"
$ws.Range("H12").Value = "'83%"
$ws.Range("H6").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("O12").Value = "10.9 °C"
$ws.Range("E13").Value = "2026-02-28 18:48:37"
$ws.Range("E14").Value = "2026-02-28 18:48:38"
$ws.Range("E15").Value = "2026-02-28 18:48:39"
$ws.Range("E16").Value = "2026-02-28 18:48:40"
$ws.Range("H16").Value = "'64%"
$ws.Range("H6").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("N16").Value = "-2.5 °C 18:00 TU"
$ws.Range("E17").Value = "2026-02-28 18:48:41"
$ws.Range("G17").Value = "2 cm"
$ws.Range("E18").Value = "2026-02-28 18:48:42"
$ws.Range("E19").Value = "2026-02-28 18:48:43"
$ws.Range("E20").Value = "2026-02-28 18:48:46"
$ws.Range("H20").Value = "'60%"
$ws.Range("H6").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("E21").Value = "2026-02-28 18:48:48"
$ws.Range("O21").Value = "7.6 °C"
$ws.Range("E22").Value = "2026-02-28 18:48:51"
$ws.Range("H22").Value = "'66%"
$ws.Range("H6").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("N22").Value = "-2.8 °C 18:26 TU"
$ws.Range("E23").Value = "2026-02-28 18:48:53"
$ws.Range("N23").Value = "-2.0 °C 18:29 TU"
$ws.Range("O23").Value = "-0.2 °C"
$ws.Range("E24").Value = "2026-02-28 18:48:55"
$ws.Range("O24").Value = "8.4 °C"
$ws.Range("E25").Value = "2026-02-28 18:48:58"
$ws.Range("H25").Value = "'59%"
$ws.Range("H6").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("N25").Value = "-0.8 °C 18:15 TU"
$ws.Range("O25").Value = "1.4 °C"
$ws.Range("E26").Value = "2026-02-28 18:49:00"
$ws.Range("H26").Value = "'79%"
$ws.Range("H6").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("E27").Value = "2026-02-28 18:49:03"
$ws.Range("H27").Value = "'52%"
$ws.Range("H6").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("N27").Value = "-0.3 °C 18:23 TU"
$ws.Range("O27").Value = "2.1 °C"
$ws.Range("E28").Value = "2026-02-28 18:49:05"
$ws.Range("H28").Value = "'81%"
$ws.Range("H6").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("O28").Value = "9.6 °C"
$ws.Range("E29").Value = "2026-02-28 18:49:07"
$ws.Range("H29").Value = "'83%"
$ws.Range("H6").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("O29").Value = "11.9 °C"
$ws.Range("E30").Value = "2026-02-28 18:49:10"
$ws.Range("E31").Value = "2026-02-28 18:49:12"
$ws.Range("H31").Value = "'82%"
$ws.Range("H6").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("E32").Value = "2026-02-28 18:49:14"
$ws.Range("E33").Value = "2026-02-28 18:49:17"
$ws.Range("H33").Value = "'64%"
$ws.Range("H6").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("E34").Value = "2026-02-28 18:49:19"
$ws.Range("H34").Value = "'65%"
$ws.Range("H6").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("O34").Value = "1.3 °C"
$ws.Range("E35").Value = "2026-02-28 18:49:22"
$ws.Range("E36").Value = "2026-02-28 18:49:24"
$ws.Range("H36").Value = "'81%"
$ws.Range("H6").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("O36").Value = "12.6 °C"
$ws.Range("E37").Value = "2026-02-28 18:49:26"
$ws.Range("O37").Value = "7.1 °C"
$ws.Range("E38").Value = "2026-02-28 18:49:29"
$ws.Range("E39").Value = "2026-02-28 18:49:31"
$ws.Range("H39").Value = "'60%"
$ws.Range("H6").Copy()
$ws.Range("H39").PasteSpecial(-4122)
$ws.Range("N39").Value = "-2.0 °C 18:04 TU"
$ws.Range("E40").Value = "2026-02-28 18:49:33"
$ws.Range("J40").Value = "1024.4 hPa"
$ws.Range("O40").Value = "7.5 °C"
$ws.Range("E41").Value = "2026-02-28 18:49:36"
$ws.Range("E42").Value = "2026-02-28 18:49:38"
$ws.Range("E43").Value = "2026-02-28 18:49:40"
$ws.Range("O43").Value = "7.5 °C"
$ws.Range("E44").Value = "2026-02-28 18:49:42"
$ws.Range("I44").Value = "1.6 mm"
$ws.Range("E45").Value = "2026-02-28 18:49:44"
$ws.Range("I45").Value = "0.7 mm"
$ws.Range("J45").Value = "1025.3 hPa"
$ws.Range("O45").Value = "6.4 °C"
$ws.Range("E46").Value = "2026-02-28 18:49:47"
$ws.Range("J46").Value = "1025.0 hPa"
$excel.CutCopyMode = $false
